$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in / clear specific "missing value" cells (imputation changes) ---

# C3: missing -> 11.2
$ws.Range("C3").Value = 11.2

# D4: -15.4 -> missing
$ws.Range("D4").Value = ""

# C5: 12.3 -> missing
$ws.Range("C5").Value = ""

# D9: missing -> -14.5
$ws.Range("D9").Value = -14.5

# D10: missing -> -14.7
$ws.Range("D10").Value = -14.7

# D17: -14.7 -> missing
$ws.Range("D17").Value = ""

# D18: -15.2 -> missing
$ws.Range("D18").Value = ""

# C21: missing -> 12.7
$ws.Range("C21").Value = 12.7

# C23: 12.2 -> missing
$ws.Range("C23").Value = ""

# --- Remove rows for "RM 232" (row 26) and "SC 92" (row 28) ---
# Delete the lower row first so the other row index stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the row deletions, the former "SC 193" row (originally row 34, now row 32)
# gets its missing C value filled in: missing -> 10.5
$ws.Range("C32").Value = 10.5
